$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.020.29'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '3.419.92'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.728'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.138'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.24'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").Value = '3.967.38'
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.141'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000212'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.81%  '
$ws.Range("D17").Value = '3.411.43'
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.48%  '
$ws.Range("D20").Value = '61.918.53'
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '489.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +20.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.77'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '34.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.88%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '12.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.167'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.72'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.89%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0492'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.44%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '150.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.78%  '
$ws.Range("E40").Value = '  +3.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.40'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.323'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.66%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.12'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.39%  '
$ws.Range("E45").Value = '  +10.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.31%  '
$ws.Range("E47").Value = '  +22.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '117.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +21.44%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.147'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +14.35%  '
